# Automatic update of files.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" (last-changed) date from 2023-09-08 (45177) to
# 2023-09-09 (45178) for every data row (rows 2-115) in column C.
for ($row = 2; $row -le 115; $row++) {
    $ws.Cells.Item($row, 3).Value = 45178
}

# Row 4 (A 13111-2022) gained a new observed species: "Smal svampklubba".
# It is a signal species, so both the signal-species count (I4) and the
# total species count (Q4) increase by one, and the species name list in
# R4 gains a new line (inserted alphabetically between "Klippfrullania"
# and "Stor revmossa").
$ws.Cells.Item(4, 9).Value = 9
$ws.Cells.Item(4, 17).Value = 13

$species = @(
    "Entita",
    "Mindre hackspett",
    "Spillkråka",
    "Talltita",
    "Blåmossa",
    "Bronshjon",
    "Fällmossa",
    "Guldlockmossa",
    "Klippfrullania",
    "Smal svampklubba",
    "Stor revmossa",
    "Västlig hakmossa",
    "Vågbandad barkbock"
)
$ws.Range("R4").Value = [string]::Join("`r`n", $species)

# Setting a new, longer wrapped value on R4 causes the row to auto-fit;
# restore the original fixed row height used throughout the sheet.
$ws.Rows.Item(4).RowHeight = 15
